$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (NumberFormat "@") for numeric-looking values so Excel
# keeps them as literal text instead of re-parsing into floats (which would
# drop things like trailing zeros, e.g. "6.90" -> 6.9, or re-render "0.0000231"
# in scientific notation).

$ws.Range("D2").Value = '61.651.06'
$ws.Range("E2").Value = '  -1.92%  '
$ws.Range("D3").Value = '2.893.13'
$ws.Range("E3").Value = '  -1.82%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.62'
$ws.Range("E5").Value = '  -4.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.87'
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("D9").Value = '2.891.53'
$ws.Range("E9").Value = '  -1.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.90'
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("E11").Value = '  -2.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.427'
$ws.Range("E12").Value = '  -2.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000231'
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.70'
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '3.373.33'
$ws.Range("E16").Value = '  -1.81%  '
$ws.Range("D17").Value = '61.639.57'
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.899.96'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.50'
$ws.Range("E19").Value = '  -2.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '430.08'
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.99'
$ws.Range("E21").Value = '  -3.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.651'
$ws.Range("E22").Value = '  -1.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.80'
$ws.Range("E23").Value = '  -2.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.98'
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.90'
$ws.Range("E27").Value = '  -11.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.00'
$ws.Range("E28").Value = '  -5.61%  '
$ws.Range("E29").Value = '  +4.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.01'
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("E31").Value = '  -4.51%  '
$ws.Range("E32").Value = '  -9.10%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.47'
$ws.Range("E35").Value = '  -3.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.957'
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.35'
$ws.Range("E37").Value = '  -4.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.78'
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("E39").Value = '  -5.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.79'
$ws.Range("E40").Value = '  -8.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.14'
$ws.Range("E41").Value = '  -3.56%  '
$ws.Range("E42").Value = '  -4.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.40'
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.266'
$ws.Range("E44").Value = '  -4.44%  '
$ws.Range("D45").Value = '2.679.00'
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.27'
$ws.Range("E46").Value = '  -2.15%  '
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '342.79'
$ws.Range("E48").Value = '  -4.13%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.45'
$ws.Range("E51").Value = '  -5.46%  '
